# Insert a new data row at row 1041, shifting existing rows 1041-1076 down
# to 1042-1077, then populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1041).Insert()

$ws.Cells.Item(1041, 1).Value = 4
$ws.Cells.Item(1041, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(1041, 3).Value = "Los Lagos"
$ws.Cells.Item(1041, 4).Value = 45075
$ws.Cells.Item(1041, 5).Value = 10
$ws.Cells.Item(1041, 6).Value = 100112020
$ws.Cells.Item(1041, 7).Value = "Tomate"
$ws.Cells.Item(1041, 8).Value = "Larga vida"
$ws.Cells.Item(1041, 9).Value = "Primera"
$ws.Cells.Item(1041, 10).Value = 250
$ws.Cells.Item(1041, 11).Value = 22000
$ws.Cells.Item(1041, 12).Value = 22000
$ws.Cells.Item(1041, 13).Value = 22000
$ws.Cells.Item(1041, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(1041, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1041, 16).Value = 1222
$ws.Cells.Item(1041, 17).Value = 18
$ws.Cells.Item(1041, 18).Value = "Hortaliza"
